$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the orphan data row (old row 5, has only F/G/H values, no A-E)
$ws.Rows.Item(5).Delete()

# Remove the old second header row (units row), now sitting at row 2
$ws.Rows.Item(2).Delete()

# Build a transient named style that mirrors the existing "Arial 9" font
# (fontId 1) but without the applyNumberFormat flag, apply it to the new
# header cells F1:K1, then drop the named style again so only the
# resulting cell-format record survives.
$headerStyle = $wb.Styles.Add("HeaderNoFmt")
$headerStyle.Font.Name = "Arial"
$headerStyle.Font.Size = 9

# Row 1 now becomes the new single header row
$ws.Range("A1:E1").ClearFormats()
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

$ws.Range("F1:K1").Style = "HeaderNoFmt"
$wb.Styles.Item("HeaderNoFmt").Delete()

# Fill in idx / idx2 columns for each power-plant row
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 108900

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 106300

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 108700

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 106400

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 106500

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 106200

$ws.Range("A4:K4").Select()
